# "cetak sepatu dan kaos"
# The document has one table with a single row and two columns, each
# holding a shoe/cap label ("T1" and "T2"). Update the label code, the
# wearer's name, the shoe/cap sizes and the class text.
#
# NOTE: Find/Replace executed against a Range obtained from a Cell's or a
# Paragraph's .Range property is not reliably scoped in this host - it can
# replace matches anywhere in the document. Re-creating the Range directly
# via $d.Range(start, end) (using the Start/End reported by the relevant
# paragraph) scopes the replacement correctly, so that pattern is used
# throughout.

$d = $word.ActiveDocument

function Replace-InParagraph($idx, $old, $new) {
    $p = $d.Paragraphs.Item($idx)
    $start = $p.Range.Start
    $end = $p.Range.End
    $scoped = $d.Range($start, $end)
    $scoped.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- First label (paragraphs 1-9): "T1" -> "C1" ---
Replace-InParagraph 1 "T1" "C1"
Replace-InParagraph 3 "SUNARTO" "ALFIAN SATYA"
Replace-InParagraph 5 "56" "57"
Replace-InParagraph 6 "DP4 NAUTIKA / 33" "DP 3 NAUTIKA/33"
Replace-InParagraph 7 "DP4 NAUTIKA / 33" "DP 3 NAUTIKA/33"

# --- Second label (paragraphs 10-17): "T2" -> "C2" ---
Replace-InParagraph 10 "T2" "C2"
Replace-InParagraph 12 "NOOR MAULANA" "RAHMAN PALA"
Replace-InParagraph 12 "41" "42"
Replace-InParagraph 14 "56" "58"
Replace-InParagraph 15 "DP4 NAUTIKA / 33" "DP 3 NAUTIKA/33"
Replace-InParagraph 16 "DP4 NAUTIKA / 33" "DP 3 NAUTIKA/33"
